# chore: update Sheets via scheduled runner
# Refreshes the cached market-price / profit figures (columns H-N:
# currentAveragePrice, currentAveragePriceNQ, currentAveragePriceHQ,
# LevePriceNQ, LevePriceHQ, LeveProfitNQ, LeveProfitHQ) for the specific
# leve rows that moved since the last scrape, across the ALC, ARM, BSM,
# CRP, CUL, GSM and LTW sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 420
$ws.Range("I28").Value = 180
$ws.Range("K28").Value = 180
$ws.Range("M28").Value = 305
# Row 33
$ws.Range("H33").Value = 3367906
$ws.Range("I33").Value = 1061.4546
$ws.Range("J33").Value = 18182022
$ws.Range("K33").Value = 1061.4546
$ws.Range("L33").Value = 18182022
$ws.Range("M33").Value = -832.4546
$ws.Range("N33").Value = -18182480
# Row 80
$ws.Range("H80").Value = 510.025
$ws.Range("I80").Value = 417.07144
$ws.Range("J80").Value = 560.0769
$ws.Range("K80").Value = 1251.21432
$ws.Range("L80").Value = 1680.2307
$ws.Range("M80").Value = -253.21432
$ws.Range("N80").Value = -3676.2307
# Row 83
$ws.Range("H83").Value = 510.025
$ws.Range("I83").Value = 417.07144
$ws.Range("J83").Value = 560.0769
$ws.Range("K83").Value = 3753.64296
$ws.Range("L83").Value = 5040.6921
$ws.Range("M83").Value = 1238.35704
$ws.Range("N83").Value = -15024.6921
# Row 86
$ws.Range("H86").Value = 1638.5
$ws.Range("I86").Value = 1435.1
$ws.Range("J86").Value = 1892.75
$ws.Range("K86").Value = 1435.1
$ws.Range("L86").Value = 1892.75
$ws.Range("M86").Value = -312.0999999999999
$ws.Range("N86").Value = -4138.75
# Row 89
$ws.Range("H89").Value = 1638.5
$ws.Range("I89").Value = 1435.1
$ws.Range("J89").Value = 1892.75
$ws.Range("K89").Value = 7175.5
$ws.Range("L89").Value = 9463.75
$ws.Range("M89").Value = -1559.5
$ws.Range("N89").Value = -20695.75
# Row 111
$ws.Range("H111").Value = 102420.1
$ws.Range("I111").Value = 2620
$ws.Range("J111").Value = 202220.2
$ws.Range("K111").Value = 7860
$ws.Range("L111").Value = 606660.6000000001
$ws.Range("M111").Value = -4793
$ws.Range("N111").Value = -612794.6000000001
# Row 129
$ws.Range("H129").Value = 804.39624
$ws.Range("J129").Value = 901.6512
$ws.Range("L129").Value = 2704.9536
$ws.Range("N129").Value = -12704.9536
# Row 137
$ws.Range("H137").Value = 1780.8636
$ws.Range("I137").Value = 1700.0834
$ws.Range("K137").Value = 5100.2502
$ws.Range("M137").Value = -2550.2502
# Row 138
$ws.Range("H138").Value = 20768.3
$ws.Range("I138").Value = 1059.7222
$ws.Range("J138").Value = 29214.834
$ws.Range("K138").Value = 3179.1666
$ws.Range("L138").Value = 87644.50199999999
$ws.Range("M138").Value = 1960.8334
$ws.Range("N138").Value = -97924.50199999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5721.9272
$ws.Range("I32").Value = 5107.282
$ws.Range("J32").Value = 7220.125
$ws.Range("K32").Value = 5107.282
$ws.Range("L32").Value = 7220.125
$ws.Range("M32").Value = -4820.282
$ws.Range("N32").Value = -7794.125
# Row 61
$ws.Range("H61").Value = 7818.8945
$ws.Range("I61").Value = 8144.9443
$ws.Range("J61").Value = 1950
$ws.Range("K61").Value = 8144.9443
$ws.Range("L61").Value = 1950
$ws.Range("M61").Value = -7932.9443
$ws.Range("N61").Value = -2374
# Row 74
$ws.Range("H74").Value = 4192.3335
$ws.Range("I74").Value = 5319.48
$ws.Range("J74").Value = 2179.5715
$ws.Range("K74").Value = 5319.48
$ws.Range("L74").Value = 2179.5715
$ws.Range("M74").Value = -4445.48
$ws.Range("N74").Value = -3927.5715
# Row 77
$ws.Range("H77").Value = 4192.3335
$ws.Range("I77").Value = 5319.48
$ws.Range("J77").Value = 2179.5715
$ws.Range("K77").Value = 26597.4
$ws.Range("L77").Value = 10897.8575
$ws.Range("M77").Value = -22229.4
$ws.Range("N77").Value = -19633.8575
# Row 88
$ws.Range("H88").Value = 83335470
$ws.Range("I88").Value = 2732.6667
$ws.Range("J88").Value = 111113060
$ws.Range("K88").Value = 2732.6667
$ws.Range("L88").Value = 111113060
$ws.Range("M88").Value = -2326.6667
$ws.Range("N88").Value = -111113872
# Row 91
$ws.Range("H91").Value = 83335470
$ws.Range("I91").Value = 2732.6667
$ws.Range("J91").Value = 111113060
$ws.Range("K91").Value = 2732.6667
$ws.Range("L91").Value = 111113060
$ws.Range("M91").Value = -1328.6667
$ws.Range("N91").Value = -111115868
# Row 132
$ws.Range("H132").Value = 4121.0347
$ws.Range("I132").Value = 1520.9487
$ws.Range("K132").Value = 4562.8461
$ws.Range("M132").Value = -2032.8461
# Row 136
$ws.Range("H136").Value = 7818.8945
$ws.Range("I136").Value = 8144.9443
$ws.Range("J136").Value = 1950
$ws.Range("K136").Value = 24434.8329
$ws.Range("L136").Value = 5850
$ws.Range("M136").Value = -21884.8329
$ws.Range("N136").Value = -10950

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 8
$ws.Range("H8").Value = 3600
$ws.Range("I8").Value = 2000
$ws.Range("J8").Value = 6000
$ws.Range("K8").Value = 2000
$ws.Range("L8").Value = 6000
$ws.Range("M8").Value = -1860
$ws.Range("N8").Value = -6280
# Row 60
$ws.Range("H60").Value = 44306
$ws.Range("J60").Value = 44306
$ws.Range("L60").Value = 44306
$ws.Range("N60").Value = -45504
# Row 86
$ws.Range("H86").Value = 2238.8333
$ws.Range("I86").Value = 2084.3333
$ws.Range("J86").Value = 2702.3333
$ws.Range("K86").Value = 2084.3333
$ws.Range("L86").Value = 2702.3333
$ws.Range("M86").Value = -961.3332999999998
$ws.Range("N86").Value = -4948.3333
# Row 89
$ws.Range("H89").Value = 2238.8333
$ws.Range("I89").Value = 2084.3333
$ws.Range("J89").Value = 2702.3333
$ws.Range("K89").Value = 10421.6665
$ws.Range("L89").Value = 13511.6665
$ws.Range("M89").Value = -4805.666499999999
$ws.Range("N89").Value = -24743.6665
# Row 134
$ws.Range("H134").Value = 3951.7307
$ws.Range("I134").Value = 4385.973
$ws.Range("K134").Value = 13157.919
$ws.Range("M134").Value = -10622.919

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 20624
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 20624
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = ""
$ws.Range("M31").Value = 20624
$ws.Range("N31").Value = -21214
# Row 34
$ws.Range("H34").Value = 20624
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 20624
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = ""
$ws.Range("M34").Value = 20624
$ws.Range("N34").Value = -21028
# Row 53
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = ""
$ws.Range("N53").Value = 0
# Row 58
$ws.Range("H58").Value = 1630.5186
$ws.Range("I58").Value = 1196.2941
$ws.Range("K58").Value = 1196.2941
$ws.Range("M58").Value = -993.2941000000001
# Row 134
$ws.Range("H134").Value = 1379.5
$ws.Range("J134").Value = 1280
$ws.Range("L134").Value = 3840
$ws.Range("N134").Value = -8910
# Row 136
$ws.Range("H136").Value = 1630.5186
$ws.Range("I136").Value = 1196.2941
$ws.Range("K136").Value = 3588.8823
$ws.Range("M136").Value = -1038.8823

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 9
$ws.Range("H9").Value = 76923850
$ws.Range("I9").Value = 400
$ws.Range("J9").Value = 90909930
$ws.Range("K9").Value = 1200
$ws.Range("L9").Value = 272729790
$ws.Range("M9").Value = -976
$ws.Range("N9").Value = -272730238
# Row 23
$ws.Range("H23").Value = 7692403
$ws.Range("I23").Value = 20000048
$ws.Range("J23").Value = 125.25
$ws.Range("K23").Value = 60000144
$ws.Range("L23").Value = 375.75
$ws.Range("M23").Value = -59999909
$ws.Range("N23").Value = -845.75
# Row 140
$ws.Range("H140").Value = 3439.55
$ws.Range("I140").Value = 3462.6843
$ws.Range("K140").Value = 10388.0529
$ws.Range("M140").Value = -5208.052899999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 4066.6897
$ws.Range("I132").Value = 4981.1
$ws.Range("J132").Value = 3585.4211
$ws.Range("K132").Value = 14943.3
$ws.Range("L132").Value = 10756.2633
$ws.Range("M132").Value = -12413.3
$ws.Range("N132").Value = -15816.2633

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 55556708
$ws.Range("I46").Value = 111111784
$ws.Range("J46").Value = 1633.3334
$ws.Range("K46").Value = 111111784
$ws.Range("L46").Value = 1633.3334
$ws.Range("M46").Value = -111111596
$ws.Range("N46").Value = -2009.3334
# Row 55
$ws.Range("H55").Value = 15625379
$ws.Range("I55").Value = 334.6154
$ws.Range("J55").Value = 26316200
$ws.Range("K55").Value = 334.6154
$ws.Range("L55").Value = 26316200
$ws.Range("M55").Value = -161.6154
$ws.Range("N55").Value = -26316546
# Row 61
$ws.Range("H61").Value = 3226.4
$ws.Range("I61").Value = 2322
$ws.Range("J61").Value = 5336.6665
$ws.Range("K61").Value = 2322
$ws.Range("L61").Value = 5336.6665
$ws.Range("M61").Value = -2120
$ws.Range("N61").Value = -5740.6665
# Row 113
$ws.Range("H113").Value = 3226.4
$ws.Range("I113").Value = 2322
$ws.Range("J113").Value = 5336.6665
$ws.Range("K113").Value = 2322
$ws.Range("L113").Value = 5336.6665
$ws.Range("M113").Value = -152
$ws.Range("N113").Value = -9676.666499999999
# Row 132
$ws.Range("H132").Value = 12042362
$ws.Range("I132").Value = 20641678
$ws.Range("J132").Value = 3318.2666
$ws.Range("K132").Value = 61925034
$ws.Range("L132").Value = 9954.799800000001
$ws.Range("M132").Value = -61922504
$ws.Range("N132").Value = -15014.7998
# Row 136
$ws.Range("H136").Value = 8737.5625
$ws.Range("I136").Value = 7286.5
$ws.Range("J136").Value = 11929.9
$ws.Range("K136").Value = 21859.5
$ws.Range("L136").Value = 35789.7
$ws.Range("M136").Value = -19309.5
$ws.Range("N136").Value = -40889.7
